# Add a new "Greece" worksheet, cloned from "Croatia", with Greece-specific
# test data, and move the "selected / active" view state from Croatia to
# the new Greece sheet (matching how Excel marks the most-recently-added
# sheet as the active tab).

$wb = $excel.ActiveWorkbook
$croatia = $wb.Worksheets.Item("Croatia")

# Clone Croatia's layout/formatting/content by copying the sheet and
# placing the copy right after it.
$croatia.Copy($null, $croatia)
$greece = $wb.Worksheets.Item($croatia.Index + 1)
$greece.Name = "Greece"

# Fill in the Greece-specific data (order matters for shared-string ids:
# NGC-4119/T3205 must be interned before "Greece Market").
$greece.Range("B4").Value = "NGC-4119/T3205"
$greece.Range("B2").Value = "Greece Market"

# Croatia's sheet view loses its "selected tab" state and its old B4
# selection becomes a full-sheet selection (as Excel does when a sheet
# stops being the active one).
$croatia.Range("A1:XFD1048576").Select()

# Greece becomes the active sheet/tab, with D14 the active cell.
$greece.Activate()
$greece.Range("D14").Select()
